$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.329.30'
$ws.Range("E2").Value = '  -3.03%  '

$ws.Range("D3").Value = '3.298.46'
$ws.Range("E3").Value = '  -3.79%  '

$ws.Range("E4").Value = '  +0.04%  '

$ws.Range("D5").Value = "'557.28"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -3.54%  '

$ws.Range("D6").Value = "'141.65"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -6.28%  '

$ws.Range("E7").Value = '  +0.06%  '

$ws.Range("D8").Value = '3.296.03'
$ws.Range("E8").Value = '  -3.82%  '

$ws.Range("E9").Value = '  -3.28%  '

$ws.Range("D10").Value = "'7.84"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -3.01%  '

$ws.Range("E11").Value = '  -4.25%  '

$ws.Range("E12").Value = '  -1.79%  '

$ws.Range("D13").Value = '3.863.29'
$ws.Range("E13").Value = '  -3.74%  '

$ws.Range("E14").Value = '  +0.18%  '

$ws.Range("D15").Value = "'26.80"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -6.09%  '

$ws.Range("D16").Value = '3.295.29'
$ws.Range("E16").Value = '  -3.53%  '

$ws.Range("E17").Value = '  -3.71%  '

$ws.Range("D18").Value = '60.329.15'
$ws.Range("E18").Value = '  -2.97%  '

$ws.Range("D19").Value = "'6.10"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -4.37%  '

$ws.Range("D20").Value = "'13.96"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -3.48%  '

$ws.Range("D21").Value = "'8.63"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -3.78%  '

$ws.Range("D22").Value = "'374.45"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.05%  '

$ws.Range("D23").Value = "'74.32"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.27%  '

$ws.Range("D24").Value = "'1.00"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.03%  '

$ws.Range("D25").Value = "'0.533"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -6.04%  '

$ws.Range("D26").Value = '3.438.49'
$ws.Range("E26").Value = '  -3.30%  '

$ws.Range("E27").Value = '  -8.79%  '

$ws.Range("D28").Value = "'0.170"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -4.84%  '

$ws.Range("D29").Value = "'1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.10%  '

$ws.Range("D30").Value = "'7.17"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -6.83%  '

$ws.Range("E31").Value = '  +0.01%  '

$ws.Range("D32").Value = "'2.03"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -4.15%  '

$ws.Range("E33").Value = '  -4.44%  '

$ws.Range("D34").Value = "'22.55"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.71%  '

$ws.Range("D35").Value = "'1.24"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -7.48%  '

$ws.Range("E36").Value = '  -6.98%  '

$ws.Range("D37").Value = "'166.78"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.93%  '

$ws.Range("D38").Value = "'1.53"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -4.11%  '

$ws.Range("D39").Value = "'6.68"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -3.26%  '

$ws.Range("B40").Value = 'RenzoRestakedETH'
$ws.Range("C40").Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range("D40").Value = '3.329.70'
$ws.Range("E40").Value = '  -3.74%  '

$ws.Range("B41").Value = 'EnergySwap'
$ws.Range("C41").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D41").Value = "'26.63"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -14.62%  '

$ws.Range("D42").Value = "'0.0727"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -6.84%  '

$ws.Range("E43").Value = '  -2.10%  '

$ws.Range("E44").Value = '  -3.63%  '

$ws.Range("E45").Value = '  -5.94%  '

$ws.Range("E46").Value = '  -5.97%  '

$ws.Range("E47").Value = '  -5.08%  '

$ws.Range("D48").Value = '2.348.69'
$ws.Range("E48").Value = '  -7.72%  '

$ws.Range("E49").Value = '  -0.01%  '

$ws.Range("D50").Value = "'6.38"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -7.31%  '

$ws.Range("D51").Value = "'21.26"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -5.95%  '
